{"js": "const pairs = [\n  [\"52\u00d712=\", \"61\u00d727=\"],\n  [\"78\u00d723=\", \"99\u00d718=\"],\n  [\"84\u00d766=\", \"68\u00d793=\"],\n  [\"31\u00d780=\", \"91\u00d729=\"],\n  [\"47\u00d726=\", \"90\u00d782=\"],\n  [\"21\u00d766=\", \"100\u00d789=\"],\n  [\"11\u00d741=\", \"31\u00d724=\"],\n  [\"24\u00d735=\", \"40\u00d736=\"],\n  [\"80\u00d789=\", \"88\u00d726=\"],\n  [\"80\u00d747=\", \"13\u00d732=\"],\n  [\"36\u00d761=\", \"33\u00d718=\"],\n  [\"63\u00d748=\", \"20\u00d723=\"],\n  [\"43\u00d782=\", \"12\u00d799=\"],\n  [\"44\u00d740=\", \"46\u00d769=\"],\n  [\"56\u00d755=\", \"13\u00d736=\"],\n  [\"17\u00d781=\", \"51\u00d766=\"],\n  [\"98\u00d794=\", \"87\u00d792=\"],\n  [\"22\u00d729=\", \"27\u00d766=\"],\n  [\"19\u00d748=\", \"45\u00d722=\"],\n  [\"22\u00d744=\", \"22\u00d778=\"],\n  [\"18\u00d724=\", \"84\u00d795=\"],\n  [\"91\u00d774=\", \"22\u00d791=\"],\n  [\"57\u00d755=\", \"58\u00d727=\"],\n  [\"67\u00d711=\", \"82\u00d771=\"],\n  [\"25\u00d725=\", \"38\u00d787=\"],\n  [\"41\u00d763=\", \"34\u00d780=\"],\n  [\"73\u00d738=\", \"94\u00d719=\"],\n  [\"50\u00d777=\", \"96\u00d737=\"],\n  [\"35\u00d720=\", \"35\u00d747=\"],\n  [\"97\u00d766=\", \"67\u00d717=\"],\n  [\"33\u00d799=\", \"74\u00d768=\"],\n  [\"27\u00d745=\", \"20\u00d766=\"],\n  [\"24\u00d712=\", \"38\u00d749=\"],\n  [\"95\u00d762=\", \"82\u00d751=\"],\n  [\"43\u00d765=\", \"100\u00d742=\"],\n  [\"27\u00d741=\", \"49\u00d786=\"],\n  [\"54\u00d737=\", \"85\u00d725=\"],\n  [\"13\u00d761=\", \"10\u00d710=\"],\n  [\"55\u00d746=\", \"65\u00d749=\"],\n  [\"67\u00d785=\", \"95\u00d781=\"],\n  [\"41\u00d789=\", \"28\u00d737=\"],\n  [\"78\u00d779=\", \"43\u00d760=\"],\n  [\"91\u00d757=\", \"22\u00d731=\"],\n  [\"25\u00d799=\", \"33\u00d797=\"],\n  [\"47\u00d769=\", \"97\u00d799=\"],\n  [\"97\u00d725=\", \"31\u00d768=\"],\n  [\"78\u00d793=\", \"51\u00d799=\"],\n  [\"71\u00d789=\", \"100\u00d757=\"],\n  [\"95\u00d744=\", \"20\u00d756=\"],\n  [\"10\u00d738=\", \"38\u00d731=\"],\n  [\"84\u00d745=\", \"74\u00d711=\"],\n  [\"85\u00d747=\", \"82\u00d733=\"],\n  [\"78\u00d736=\", \"19\u00d717=\"],\n  [\"93\u00d789=\", \"29\u00d782=\"],\n  [\"43\u00d757=\", \"26\u00d723=\"],\n  [\"38\u00d748=\", \"42\u00d795=\"],\n  [\"71\u00d726=\", \"81\u00d772=\"],\n  [\"53\u00d757=\", \"95\u00d753=\"],\n  [\"32\u00d744=\", \"96\u00d710=\"],\n  [\"51\u00d769=\", \"15\u00d770=\"],\n  [\"57\u00d747=\", \"80\u00d728=\"],\n  [\"56\u00d765=\", \"54\u00d756=\"],\n  [\"93\u00d767=\", \"82\u00d771=\"],\n  [\"72\u00d748=\", \"54\u00d719=\"],\n  [\"68\u00d758=\", \"55\u00d715=\"],\n  [\"14\u00d797=\", \"55\u00d760=\"],\n  [\"17\u00d720=\", \"50\u00d755=\"],\n  [\"35\u00d737=\", \"68\u00d775=\"],\n  [\"69\u00d733=\", \"57\u00d781=\"],\n  [\"33\u00d777=\", \"66\u00d767=\"],\n  [\"90\u00d799=\", \"16\u00d794=\"],\n  [\"65\u00d768=\", \"83\u00d785=\"],\n  [\"78\u00d762=\", \"51\u00d742=\"],\n  [\"44\u00d786=\", \"40\u00d789=\"],\n  [\"11\u00d724=\", \"25\u00d721=\"],\n  [\"50\u00d730=\", \"33\u00d714=\"],\n  [\"60\u00d738=\", \"30\u00d796=\"],\n  [\"68\u00d722=\", \"90\u00d7100=\"],\n  [\"96\u00d733=\", \"71\u00d767=\"],\n  [\"13\u00d729=\", \"79\u00d745=\"],\n  [\"22\u00d724=\", \"28\u00d777=\"],\n  [\"53\u00d737=\", \"19\u00d741=\"],\n  [\"38\u00d737=\", \"38\u00d755=\"],\n  [\"95\u00d782=\", \"30\u00d763=\"],\n  [\"76\u00d762=\", \"41\u00d7100=\"],\n  [\"78\u00d713=\", \"59\u00d780=\"],\n  [\"73\u00d789=\", \"71\u00d767=\"],\n  [\"60\u00d793=\", \"30\u00d761=\"],\n  [\"12\u00d777=\", \"54\u00d793=\"],\n  [\"37\u00d741=\", \"36\u00d789=\"],\n  [\"88\u00d785=\", \"61\u00d799=\"],\n  [\"28\u00d794=\", \"98\u00d713=\"],\n  [\"64\u00d797=\", \"37\u00d738=\"],\n  [\"72\u00d719=\", \"38\u00d780=\"],\n  [\"41\u00d743=\", \"86\u00d720=\"],\n  [\"50\u00d786=\", \"47\u00d795=\"],\n  [\"50\u00d793=\", \"67\u00d790=\"],\n  [\"19\u00d710=\", \"82\u00d757=\"],\n  [\"84\u00d791=\", \"36\u00d715=\"],\n  [\"89\u00d797=\", \"41\u00d777=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"52\u00d712=\", \"61\u00d727=\"),\n    @(\"78\u00d723=\", \"99\u00d718=\"),\n    @(\"84\u00d766=\", \"68\u00d793=\"),\n    @(\"31\u00d780=\", \"91\u00d729=\"),\n    @(\"47\u00d726=\", \"90\u00d782=\"),\n    @(\"21\u00d766=\", \"100\u00d789=\"),\n    @(\"11\u00d741=\", \"31\u00d724=\"),\n    @(\"24\u00d735=\", \"40\u00d736=\"),\n    @(\"80\u00d789=\", \"88\u00d726=\"),\n    @(\"80\u00d747=\", \"13\u00d732=\"),\n    @(\"36\u00d761=\", \"33\u00d718=\"),\n    @(\"63\u00d748=\", \"20\u00d723=\"),\n    @(\"43\u00d782=\", \"12\u00d799=\"),\n    @(\"44\u00d740=\", \"46\u00d769=\"),\n    @(\"56\u00d755=\", \"13\u00d736=\"),\n    @(\"17\u00d781=\", \"51\u00d766=\"),\n    @(\"98\u00d794=\", \"87\u00d792=\"),\n    @(\"22\u00d729=\", \"27\u00d766=\"),\n    @(\"19\u00d748=\", \"45\u00d722=\"),\n    @(\"22\u00d744=\", \"22\u00d778=\"),\n    @(\"18\u00d724=\", \"84\u00d795=\"),\n    @(\"91\u00d774=\", \"22\u00d791=\"),\n    @(\"57\u00d755=\", \"58\u00d727=\"),\n    @(\"67\u00d711=\", \"82\u00d771=\"),\n    @(\"25\u00d725=\", \"38\u00d787=\"),\n    @(\"41\u00d763=\", \"34\u00d780=\"),\n    @(\"73\u00d738=\", \"94\u00d719=\"),\n    @(\"50\u00d777=\", \"96\u00d737=\"),\n    @(\"35\u00d720=\", \"35\u00d747=\"),\n    @(\"97\u00d766=\", \"67\u00d717=\"),\n    @(\"33\u00d799=\", \"74\u00d768=\"),\n    @(\"27\u00d745=\", \"20\u00d766=\"),\n    @(\"24\u00d712=\", \"38\u00d749=\"),\n    @(\"95\u00d762=\", \"82\u00d751=\"),\n    @(\"43\u00d765=\", \"100\u00d742=\"),\n    @(\"27\u00d741=\", \"49\u00d786=\"),\n    @(\"54\u00d737=\", \"85\u00d725=\"),\n    @(\"13\u00d761=\", \"10\u00d710=\"),\n    @(\"55\u00d746=\", \"65\u00d749=\"),\n    @(\"67\u00d785=\", \"95\u00d781=\"),\n    @(\"41\u00d789=\", \"28\u00d737=\"),\n    @(\"78\u00d779=\", \"43\u00d760=\"),\n    @(\"91\u00d757=\", \"22\u00d731=\"),\n    @(\"25\u00d799=\", \"33\u00d797=\"),\n    @(\"47\u00d769=\", \"97\u00d799=\"),\n    @(\"97\u00d725=\", \"31\u00d768=\"),\n    @(\"78\u00d793=\", \"51\u00d799=\"),\n    @(\"71\u00d789=\", \"100\u00d757=\"),\n    @(\"95\u00d744=\", \"20\u00d756=\"),\n    @(\"10\u00d738=\", \"38\u00d731=\"),\n    @(\"84\u00d745=\", \"74\u00d711=\"),\n    @(\"85\u00d747=\", \"82\u00d733=\"),\n    @(\"78\u00d736=\", \"19\u00d717=\"),\n    @(\"93\u00d789=\", \"29\u00d782=\"),\n    @(\"43\u00d757=\", \"26\u00d723=\"),\n    @(\"38\u00d748=\", \"42\u00d795=\"),\n    @(\"71\u00d726=\", \"81\u00d772=\"),\n    @(\"53\u00d757=\", \"95\u00d753=\"),\n    @(\"32\u00d744=\", \"96\u00d710=\"),\n    @(\"51\u00d769=\", \"15\u00d770=\"),\n    @(\"57\u00d747=\", \"80\u00d728=\"),\n    @(\"56\u00d765=\", \"54\u00d756=\"),\n    @(\"93\u00d767=\", \"82\u00d771=\"),\n    @(\"72\u00d748=\", \"54\u00d719=\"),\n    @(\"68\u00d758=\", \"55\u00d715=\"),\n    @(\"14\u00d797=\", \"55\u00d760=\"),\n    @(\"17\u00d720=\", \"50\u00d755=\"),\n    @(\"35\u00d737=\", \"68\u00d775=\"),\n    @(\"69\u00d733=\", \"57\u00d781=\"),\n    @(\"33\u00d777=\", \"66\u00d767=\"),\n    @(\"90\u00d799=\", \"16\u00d794=\"),\n    @(\"65\u00d768=\", \"83\u00d785=\"),\n    @(\"78\u00d762=\", \"51\u00d742=\"),\n    @(\"44\u00d786=\", \"40\u00d789=\"),\n    @(\"11\u00d724=\", \"25\u00d721=\"),\n    @(\"50\u00d730=\", \"33\u00d714=\"),\n    @(\"60\u00d738=\", \"30\u00d796=\"),\n    @(\"68\u00d722=\", \"90\u00d7100=\"),\n    @(\"96\u00d733=\", \"71\u00d767=\"),\n    @(\"13\u00d729=\", \"79\u00d745=\"),\n    @(\"22\u00d724=\", \"28\u00d777=\"),\n    @(\"53\u00d737=\", \"19\u00d741=\"),\n    @(\"38\u00d737=\", \"38\u00d755=\"),\n    @(\"95\u00d782=\", \"30\u00d763=\"),\n    @(\"76\u00d762=\", \"41\u00d7100=\"),\n    @(\"78\u00d713=\", \"59\u00d780=\"),\n    @(\"73\u00d789=\", \"71\u00d767=\"),\n    @(\"60\u00d793=\", \"30\u00d761=\"),\n    @(\"12\u00d777=\", \"54\u00d793=\"),\n    @(\"37\u00d741=\", \"36\u00d789=\"),\n    @(\"88\u00d785=\", \"61\u00d799=\"),\n    @(\"28\u00d794=\", \"98\u00d713=\"),\n    @(\"64\u00d797=\", \"37\u00d738=\"),\n    @(\"72\u00d719=\", \"38\u00d780=\"),\n    @(\"41\u00d743=\", \"86\u00d720=\"),\n    @(\"50\u00d786=\", \"47\u00d795=\"),\n    @(\"50\u00d793=\", \"67\u00d790=\"),\n    @(\"19\u00d710=\", \"82\u00d757=\"),\n    @(\"84\u00d791=\", \"36\u00d715=\"),\n    @(\"89\u00d797=\", \"41\u00d777=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $found = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
